# Grado04.xlsx / seguimiento ("Seguimiento" sheet) update:
#   - Fill in the tracking dates for row 13 (item 8, CN_04_08_CO):
#       B13 (Manuscrito)            -> 2015-03-10
#       C13 (Publicación manuscrito) -> 2015-03-12
#       D13 (Formatos de recursos)   -> 2015-03-21
#   - Leave the cursor on E13, matching where the editor left off.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B13").Value = 42073
$ws.Range("C13").Value = 42075
$ws.Range("D13").Value = 42084

$ws.Range("E13").Select() | Out-Null
